# Config UART2 with UART1
# Updates the GPIO pin-assignment bullet list (numId=4) so that UART1 and
# UART2 each get a dedicated Tx pin, Rx pin, and "response of received
# data" pin, instead of the old mixed-up photo-register / computer wording.

$d = $word.ActiveDocument

function Get-ParagraphByText($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd("`r")
        if ($t -eq $text) {
            return $p
        }
    }
    throw "Paragraph with text '$text' not found"
}

function Replace-Paragraph-Text($oldText, $newText) {
    $p = Get-ParagraphByText $oldText
    $start = $p.Range.Start
    $rng = $d.Range($start, $start + $oldText.Length)
    $rng.Text = $newText
}

# 1) PA0: drop the "/RX" -- it's now TX-only (RX gets its own bullet below)
Replace-Paragraph-Text "PA0: TX/RX signal for uart1 : GPIO output pushpull" `
                        "PA0: TX signal for uart1 : GPIO output pushpull"

# 2) PA1: used to be the "photo reg" bullet -- repurpose it as the UART1 Rx
#    bullet. Keep the "PA1: " prefix and set the rest to the new text. The
#    description is replaced first so the prefix's offsets stay valid.
$oldPa1 = "PA1: Output in response of photo reg : GPIO output pushpull"
$p1 = Get-ParagraphByText $oldPa1
$p1Start = $p1.Range.Start
$prefixLen = "PA1: ".Length
$descRng = $d.Range($p1Start + $prefixLen, $p1Start + $oldPa1.Length)
$descRng.Text = "RX signal for uart1 : GPIO output pushpull"

# 3) Old PA5 bullet becomes the new PA4 bullet (response of received data
#    from USART1), trailing space included as in the source.
Replace-Paragraph-Text "PA5: Tx/Rx signal for UART2: GPIO output push pull" `
                        "PA4: Response of received data from USART1 : GPIO OUTput pushpull "

# 4) Old pA6 bullet becomes the PA5 / UART2 Tx bullet.
Replace-Paragraph-Text "pA6: Output is response of computer: GPIO output pushpull" `
                        "PA5: TxRx signal for UART2: GPIO output push pull"

# 5) Old PA7 bullet (2 runs: "PA7: Input for photo register" + ": analog
#    input") becomes the new pA6 / UART2 Rx bullet. Replace the
#    description first (keeps the "PA7: " prefix range's offsets valid),
#    then overwrite the prefix text itself.
$oldPa7 = "PA7: Input for photo register: analog input"
$p7 = Get-ParagraphByText $oldPa7
$p7Start = $p7.Range.Start
$oldPrefixLen = "PA7: ".Length
$descRng7 = $d.Range($p7Start + $oldPrefixLen, $p7Start + $oldPa7.Length)
$descRng7.Text = "Rx signal for UART2: GPIO output push pull"
$prefixRng7 = $d.Range($p7Start, $p7Start + $oldPrefixLen)
$prefixRng7.Text = "pA6: "

# 6) Append a brand-new bullet for PA7 (response of received data from
#    USART2) right after the bullet we just rewrote.
$p7.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($p7.Index + 1)
$newPara.Range.Text = "PA7: Response of received data from USART2 : GPIO output pushpull"

Write-Output "done"
